# Regenerate merged AHB files
#
# This mirrors the commit that re-ran the AHB-diff generator against the
# FV2310 / FV2404 message-format pairing instead of a generic "old"/"new"
# pairing:
#   - every "<Label>_old" header becomes "<Label>_FV2310"
#   - every "<Label>_new" header becomes "<Label>_FV2404"
#   - the data range gets wrapped in a native Excel Table ("Table1")
#   - the header row is frozen so it stays visible while scrolling

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row (old -> FV2310, new -> FV2404) --------------
$headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# --- 2. Turn the used range into a real Table (Table1), A1:U66 -------------
$dataRange = $ws.Range("A1:U66")
$tbl = $ws.ListObjects.Add(1, $dataRange, $false, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row (split below row 1, top-left cell A2) -------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "Renamed headers, added Table1 over $($dataRange.Address()), froze header row."
